$d = $word.ActiveDocument

$replacements = @(
    @("2025-02-12 Wednesday", "2025-02-13 Thursday"),
    @("63×34=2142", "85×23=1955"),
    @("32×62=1984", "39×98=3822"),
    @("34×35=1190", "33×84=2772"),
    @("81×35=2835", "94×92=8648"),
    @("56×55=3080", "75×65=4875"),
    @("30×57=1710", "92×65=5980"),
    @("15×41=615", "12×14=168"),
    @("73×49=3577", "54×34=1836"),
    @("46×95=4370", "56×52=2912"),
    @("14×33=462", "20×14=280"),
    @("24×73=1752", "88×14=1232"),
    @("69×40=2760", "41×70=2870"),
    @("86×44=3784", "89×42=3738"),
    @("94×19=1786", "22×55=1210"),
    @("26×84=2184", "78×72=5616"),
    @("69×53=3657", "84×62=5208"),
    @("96×27=2592", "13×89=1157"),
    @("83×36=2988", "84×57=4788"),
    @("12×54=648", "20×35=700"),
    @("36×11=396", "43×33=1419"),
    @("96×40=3840", "91×65=5915"),
    @("79×73=5767", "16×40=640"),
    @("45×79=3555", "56×26=1456"),
    @("39×83=3237", "28×49=1372"),
    @("29×76=2204", "89×57=5073")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
